# Update the document date header, then the 25 two-digit-by-two-digit
# multiplication problems scattered across the tables. Every "old" text
# below is unique in the document, so a simple global Find/Replace
# (wdReplaceAll = 2) on $d.Content is sufficient and safe.
$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-08-31 Saturday", $false, $false, $false, $false, $false, $true, 1, $false, "2024-09-01 Sunday", 2) | Out-Null
$d.Content.Find.Execute("17×54=918", $false, $false, $false, $false, $false, $true, 1, $false, "44×17=748", 2) | Out-Null
$d.Content.Find.Execute("91×92=8372", $false, $false, $false, $false, $false, $true, 1, $false, "11×94=1034", 2) | Out-Null
$d.Content.Find.Execute("33×56=1848", $false, $false, $false, $false, $false, $true, 1, $false, "72×48=3456", 2) | Out-Null
$d.Content.Find.Execute("29×68=1972", $false, $false, $false, $false, $false, $true, 1, $false, "12×36=432", 2) | Out-Null
$d.Content.Find.Execute("16×39=624", $false, $false, $false, $false, $false, $true, 1, $false, "86×34=2924", 2) | Out-Null
$d.Content.Find.Execute("46×37=1702", $false, $false, $false, $false, $false, $true, 1, $false, "52×68=3536", 2) | Out-Null
$d.Content.Find.Execute("88×90=7920", $false, $false, $false, $false, $false, $true, 1, $false, "12×66=792", 2) | Out-Null
$d.Content.Find.Execute("48×63=3024", $false, $false, $false, $false, $false, $true, 1, $false, "12×48=576", 2) | Out-Null
$d.Content.Find.Execute("81×61=4941", $false, $false, $false, $false, $false, $true, 1, $false, "82×23=1886", 2) | Out-Null
$d.Content.Find.Execute("78×78=6084", $false, $false, $false, $false, $false, $true, 1, $false, "31×96=2976", 2) | Out-Null
$d.Content.Find.Execute("96×44=4224", $false, $false, $false, $false, $false, $true, 1, $false, "56×29=1624", 2) | Out-Null
$d.Content.Find.Execute("50×99=4950", $false, $false, $false, $false, $false, $true, 1, $false, "72×64=4608", 2) | Out-Null
$d.Content.Find.Execute("20×89=1780", $false, $false, $false, $false, $false, $true, 1, $false, "46×59=2714", 2) | Out-Null
$d.Content.Find.Execute("67×84=5628", $false, $false, $false, $false, $false, $true, 1, $false, "80×60=4800", 2) | Out-Null
$d.Content.Find.Execute("81×55=4455", $false, $false, $false, $false, $false, $true, 1, $false, "95×89=8455", 2) | Out-Null
$d.Content.Find.Execute("39×96=3744", $false, $false, $false, $false, $false, $true, 1, $false, "74×80=5920", 2) | Out-Null
$d.Content.Find.Execute("59×27=1593", $false, $false, $false, $false, $false, $true, 1, $false, "65×20=1300", 2) | Out-Null
$d.Content.Find.Execute("61×17=1037", $false, $false, $false, $false, $false, $true, 1, $false, "74×65=4810", 2) | Out-Null
$d.Content.Find.Execute("22×89=1958", $false, $false, $false, $false, $false, $true, 1, $false, "24×49=1176", 2) | Out-Null
$d.Content.Find.Execute("54×37=1998", $false, $false, $false, $false, $false, $true, 1, $false, "27×21=567", 2) | Out-Null
$d.Content.Find.Execute("17×72=1224", $false, $false, $false, $false, $false, $true, 1, $false, "78×45=3510", 2) | Out-Null
$d.Content.Find.Execute("91×63=5733", $false, $false, $false, $false, $false, $true, 1, $false, "17×36=612", 2) | Out-Null
$d.Content.Find.Execute("36×64=2304", $false, $false, $false, $false, $false, $true, 1, $false, "81×47=3807", 2) | Out-Null
$d.Content.Find.Execute("16×40=640", $false, $false, $false, $false, $false, $true, 1, $false, "86×17=1462", 2) | Out-Null
$d.Content.Find.Execute("40×94=3760", $false, $false, $false, $false, $false, $true, 1, $false, "23×51=1173", 2) | Out-Null
